# Updating the GESS model with the latest data
#
# The workbook's "Lookup" helper column (shared strings) contains date
# labels for two days; those two days move forward by 5 days
# (20/21 Feb 2026 -> 25/26 Feb 2026). Column A (Timestamp) shifts by the
# same 5 days, and column B (Actual Consumption (MW)) is refreshed with
# newly fetched values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "Lookup" text labels stored in the shared strings table ---
# These are plain text cells in column D, formatted like "20.02.2026<quarter>".
# Replace the old day labels with the new ones (order matters so the
# already-replaced "25.02.2026" text from day 1 isn't matched again while
# processing day 2).
$ws.Cells.Replace("20.02.2026", "25.02.2026")
$ws.Cells.Replace("21.02.2026", "26.02.2026")

# --- 2. Shift the Timestamp column (A2:A193) forward by 5 days ---
$rngA = $ws.Range("A2:A193")
$valsA = $rngA.Value2
for ($i = 1; $i -le $valsA.GetUpperBound(0); $i++) {
    $valsA[$i, 1] = $valsA[$i, 1] + 5
}
$rngA.Value = $valsA

# --- 3. Refresh the Actual Consumption (MW) column (B2:B193) with new data ---
$newConsumption = @(
    6070,6041,5998,5999,6073,5980,5943,5888,5894,5899,5950,5971,5941,5933,6005,
    6026,6131,6158,6277,6376,6555,6658,6837,6962,7193,7327,7437,7480,7614,7605,
    7557,7397,7434,7477,7323,7292,7147,7100,6977,6913,6857,6810,6775,6749,6676,
    6704,6613,6683,6696,6797,6801,6803,6947,6972,6992,6951,6961,6951,6991,7029,
    7163,7263,7370,7406,7543,7614,7749,7847,7961,8106,8203,0,8220,8200,8215,
    8143,8120,8077,8005,7911,7788,7694,7562,7365,7199,7015,6908,6705,6598,6522,
    6430,0,6262,6217,6168,6113,6125,6048,6031,0,5978,5929,5912,5898,5886,5895,
    5856,5925,5924,5955,5935,5960,6087,6137,6226,6328,6546,6707,6791,6971,7154,
    7393,7498,7558,7633,7648,7607,7569,7562,7518,7459,7365,7149,7015,6910,0,0,
    0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,
    0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0
)

$rngB = $ws.Range("B2:B193")
$valsB = $rngB.Value2
for ($i = 1; $i -le $valsB.GetUpperBound(0); $i++) {
    $valsB[$i, 1] = $newConsumption[$i - 1]
}
$rngB.Value = $valsB
